# Fruta / hortaliza, semanal
# Insert a new weekly record row at row 236 (pushing existing rows 236:273 down
# to 237:274), carrying over the same dimension/category metadata as the
# (former) row 236, but with this week's own price/volume figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 236:273 down to 237:274, leaving a blank row 236 to populate.
$ws.Rows("236:236").Insert()

# Populate the new row 236 with the new weekly observation.
$ws.Cells.Item(236, 1).Value = 9
$ws.Cells.Item(236, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(236, 3).Value = "Metropolitana"
$ws.Cells.Item(236, 4).Value = 44491
$ws.Cells.Item(236, 5).Value = 13
$ws.Cells.Item(236, 6).Value = "Fruta"
$ws.Cells.Item(236, 7).Value = 100108
$ws.Cells.Item(236, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(236, 9).Value = 100108002
$ws.Cells.Item(236, 10).Value = "Mango"
$ws.Cells.Item(236, 11).Value = "Sin especificar"
$ws.Cells.Item(236, 12).Value = "Primera"
$ws.Cells.Item(236, 13).Value = 450
$ws.Cells.Item(236, 14).Value = 6000
$ws.Cells.Item(236, 15).Value = 6000
$ws.Cells.Item(236, 16).Value = 6000
$ws.Cells.Item(236, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(236, 18).Value = "Perú"
$ws.Cells.Item(236, 19).Value = 1500
$ws.Cells.Item(236, 20).Value = 4

# Match the date formatting/style used by the other "Fecha" column cells.
$ws.Cells.Item(236, 4).NumberFormat = $ws.Cells.Item(237, 4).NumberFormat
